# Atualização de bases das ligas, do dia: 17-02-2024 às 22:47
#
# Several rows had their match-detail data (columns B:AC) swapped between
# two adjacent rows while the running index in column A stayed put.
# Row pairs affected: (14,15) (24,25) (26,27) (31,32) (132,133)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(14, 15),
    @(24, 25),
    @(26, 27),
    @(31, 32),
    @(132, 133)
)

foreach ($pair in $pairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    $addr1 = "B" + $row1 + ":AC" + $row1
    $addr2 = "B" + $row2 + ":AC" + $row2

    $range1 = $ws.Range($addr1)
    $range2 = $ws.Range($addr2)

    # Capture current values before overwriting anything.
    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value = $vals2
    $range2.Value = $vals1
}
